# Final Answer for commits
#
# The paragraph that used to read:
#   "We need to stage commit and push after reverting because we changed
#    on our local machine but on github nothing changed so we need to
#    publish the modifications by pushing the new file as soon as we
#    changed."
# is duplicated as a brand-new paragraph placed right before itself, and
# the original paragraph's text is replaced with "I have 11 commits."
# (keeping the original paragraph's identity/bookmark and the formatting
# of its first run).

$d = $word.ActiveDocument

# Locate the paragraph holding the target sentence (robust to its
# position in the document instead of assuming it is the last one).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*We need to stage commit and push after*") {
        $target = $p
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found - aborting."
} else {
    $range = $target.Range

    # Brand-new paragraph (no rsid/paraId - it did not exist before),
    # duplicating the original wording/formatting exactly.
    $newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="BDC1C6"/><w:shd w:val="clear" w:color="auto" w:fill="202124"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="BDC1C6"/><w:shd w:val="clear" w:color="auto" w:fill="202124"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">We need to stage commit and push after </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="BDC1C6"/><w:shd w:val="clear" w:color="auto" w:fill="202124"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">reverting because we changed on our local machine but on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="BDC1C6"/><w:shd w:val="clear" w:color="auto" w:fill="202124"/><w:lang w:val="en-US"/></w:rPr><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="BDC1C6"/><w:shd w:val="clear" w:color="auto" w:fill="202124"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> nothing changed so we need to publish the modifications by pushing the new file</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="BDC1C6"/><w:shd w:val="clear" w:color="auto" w:fill="202124"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> as soon as we changed.</w:t></w:r></w:p>
'@

    # The original paragraph, keeping its paragraph mark/bookmark, its
    # pPr, and the rPr of its first run, but with the new wording.
    $updatedOriginalParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="65EAF851" w14:textId="0B978C6B" w:rsidR="008578F2" w:rsidRPr="009149BF" w:rsidRDefault="008578F2"><w:pPr><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="BDC1C6"/><w:shd w:val="clear" w:color="auto" w:fill="202124"/><w:lang w:val="en-US"/></w:rPr><w:t>I have 11 commits.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

    $range.InsertXML($newParagraphXml + $updatedOriginalParagraphXml)
}
